# Daily attendance processing - 2026-01-25 14:59:16
# Swap the order of "System" and the email address in the "Recorded By"
# column (column G) from "System, dnasr281@gmail.com" to
# "dnasr281@gmail.com, System" for every row where that exact text appears.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$oldText = "System, dnasr281@gmail.com"
$newText = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    if ($cell.Value2 -eq $oldText) {
        $cell.Value = $newText
    }
}
